# Rule "R30" (row 10) had its "Integer min" value (column C) changed
# from 18 to 1 on the "Rules" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
